# Corrected an english error
#
# 1) Bump the cached "datetimeFigureOut" footer-date field from 05/12/2023
#    to 06/12/2023 everywhere it is cached: every slide layout's date
#    placeholder, and the notes master's date placeholder.
# 2) Fix two small English wording issues on slide 6's "Requirements:"
#    text box:
#       "filtering the ingredients"          -> "filtering by the ingredients"
#       "filtering with different settings"  -> "filtering by different settings"

$p = $ppt.ActivePresentation

$oldDate = "05/12/2023"
$newDate = "06/12/2023"

# --- Update the date placeholder cached text on every slide layout ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Update the date placeholder cached text on the notes master ---
# (Going through Shapes() directly is unreliable for the notes master in
# this host, so use the HeadersFooters facade which targets the correct
# part.)
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDate

# --- Fix the wording on slide 6 ("Requirements and Entities ... ") ---
$slide6 = $p.Slides.Item(6)
$reqBox = $slide6.Shapes.Item(3)
$tr = $reqBox.TextFrame.TextRange

for ($k = 1; $k -le $tr.Paragraphs().Count; $k++) {
    $para = $tr.Paragraphs($k, 1)
    $ptext = $para.Text.TrimEnd("`r")
    if ($ptext -eq "Allow users to look for new recipes filtering the ingredients.") {
        $para.Runs(1, 1).Text = "Allow users to look for new recipes filtering by the ingredients."
    }
    elseif ($ptext -eq "Find users filtering with different settings defined by the admin.") {
        $para.Runs(1, 1).Text = "Find users filtering by different settings defined by the admin."
    }
}
